# feat: add 2022-Q4 data
#
# Inserts a new worksheet "2022-Q4" right after "总计" (holding the Q4 fund
# holdings detail, mirroring the layout of the existing "2022-Q3" /
# "2022-Q1" sheets) and updates the "总计" summary sheet with the new
# 2022-Q4 totals row (pushing the existing 2022-Q3 / 2022-Q1 rows down).

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet directly after "总计" (so the tab
#    order becomes 总计, 2022-Q4, 2022-Q3, 2022-Q1).
# ---------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Add([Type]::Missing, $wsTotal)
$wsQ4.Name = "2022-Q4"

# Re-fetch by name (not captured before the Add() above) since inserting
# a sheet shifts the Worksheets collection and stale references resolve
# to the wrong sheet afterwards.
$wsQ3 = $wb.Worksheets.Item("2022-Q3")   # used purely as a style donor

# Match the page margins used by the sibling sheets (0.75/0.75/1/1/.5/.5 in).
$wsQ4.PageSetup.LeftMargin = 54
$wsQ4.PageSetup.RightMargin = 54
$wsQ4.PageSetup.TopMargin = 72
$wsQ4.PageSetup.BottomMargin = 72
$wsQ4.PageSetup.HeaderMargin = 36
$wsQ4.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsQ4.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 3. Data rows (code / name / size / position / ratio / value / rank).
#    Columns B-G are kept as plain text (as in the sibling sheets) so
#    fund codes like "001227" don't get coerced into numbers.
# ---------------------------------------------------------------------
$rows = @(
    @("001227", "中邮信息产业灵活配置混合", "5.96", "84.02", "3.71", "0.2211", 5),
    @("011346", "淳厚鑫淳一年持有期混合", "3.34", "78.96", "3.05", "0.1019", 5),
    @("012454", "淳厚鑫悦混合A", "1.82", "85.29", "3.91", "0.0712", 2),
    @("001275", "中邮创新优势灵活配置混合", "1.06", "83.76", "3.63", "0.0385", 6),
    @("002145", "诺安景鑫灵活配置混合", "0.50", "83.31", "5.14", "0.0257", 4),
    @("012455", "淳厚鑫悦混合C", "0.57", "85.29", "3.91", "0.0223", 2),
    @("003308", "中信建投睿利灵活配置混合A", "0.07", "71.26", "4.75", "0.0033", 3),
    @("004635", "中信建投睿利灵活配置混合C", "0.04", "71.26", "4.75", "0.0019", 3)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $r + 2
    $wsQ4.Cells.Item($row, 1).Value = $r

    $wsQ4.Cells.Item($row, 2).NumberFormat = "@"
    $wsQ4.Cells.Item($row, 2).Value = $rows[$r][0]

    $wsQ4.Cells.Item($row, 3).Value = $rows[$r][1]

    $wsQ4.Cells.Item($row, 4).NumberFormat = "@"
    $wsQ4.Cells.Item($row, 4).Value = $rows[$r][2]

    $wsQ4.Cells.Item($row, 5).NumberFormat = "@"
    $wsQ4.Cells.Item($row, 5).Value = $rows[$r][3]

    $wsQ4.Cells.Item($row, 6).NumberFormat = "@"
    $wsQ4.Cells.Item($row, 6).Value = $rows[$r][4]

    $wsQ4.Cells.Item($row, 7).NumberFormat = "@"
    $wsQ4.Cells.Item($row, 7).Value = $rows[$r][5]

    $wsQ4.Cells.Item($row, 8).Value = $rows[$r][6]
}

# ---------------------------------------------------------------------
# 4. Formatting: reuse the bold / bordered / centred style that already
#    marks the header row + first (index) column on the sibling sheets,
#    via a plain format copy/paste so no duplicate style entries are
#    created.
# ---------------------------------------------------------------------
$wsQ3.Range("B1:H1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)

$wsQ3.Range("A2").Copy()
$wsQ4.Range("A2:A9").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5. Update the "总计" overview sheet: shift the existing 2022-Q3 / 2022-Q1
#    rows down by one and insert the new 2022-Q4 totals as row 2.
# ---------------------------------------------------------------------

# Row 4 <- old row 3 (2022-Q1), index bumped from 1 to 2.
$wsTotal.Cells.Item(3, 1).Copy()
$wsTotal.Cells.Item(4, 1).PasteSpecial(-4122)
$wsTotal.Cells.Item(4, 1).Value = 2
$wsTotal.Cells.Item(4, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(4, 3).Value = 2
$wsTotal.Cells.Item(4, 4).Value = 0

# Row 3 <- old row 2 (2022-Q3), index bumped from 0 to 1.
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(3, 3).Value = 6
$wsTotal.Cells.Item(3, 4).Value = 0.04

# Row 2 <- new 2022-Q4 totals.
$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"
$wsTotal.Cells.Item(2, 3).Value = 8
$wsTotal.Cells.Item(2, 4).Value = 0.49

# ---------------------------------------------------------------------
# 6. Restore "2022-Q1" as the selected/active tab (it was the active
#    sheet before this edit and the edit doesn't intend to change that).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
